$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '63.372.24'
$ws.Range('E2').Value2 = '  +3.94%  '
$ws.Range('D3').Value2 = '3.487.06'
$ws.Range('E3').Value2 = '  +3.36%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value2 = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '584.29'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value2 = '  +2.00%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value2 = '147.90'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value2 = '  +6.60%  '
$ws.Range('E7').Value2 = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.478'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value2 = '  +1.16%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '7.72'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value2 = '  +0.48%  '
$ws.Range('E10').Value2 = '  +4.01%  '
$ws.Range('E11').Value2 = '  +3.64%  '
$ws.Range('D12').Value2 = '4.082.25'
$ws.Range('E12').Value2 = '  +3.37%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '29.73'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value2 = '  +5.76%  '
$ws.Range('E14').Value2 = '  -0.41%  '
$ws.Range('D15').Value2 = '3.523.12'
$ws.Range('E15').Value2 = '  +4.68%  '
$ws.Range('E16').Value2 = '  +3.27%  '
$ws.Range('D17').Value2 = '63.337.21'
$ws.Range('E17').Value2 = '  +3.78%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '6.28'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value2 = '  +3.16%  '
$ws.Range('E19').Value2 = '  +6.18%  '
$ws.Range('E20').Value2 = '  +4.93%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '391.38'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value2 = '  +1.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '0.564'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value2 = '  +2.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '75.23'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value2 = '  +0.35%  '
$ws.Range('E24').Value2 = '  -0.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '0.0000119'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value2 = '  +7.22%  '
$ws.Range('D26').Value2 = '3.627.77'
$ws.Range('E26').Value2 = '  +3.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '0.183'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value2 = '  -2.98%  '
$ws.Range('E28').Value2 = '  +9.74%  '
$ws.Range('E29').Value2 = '  -0.05%  '
$ws.Range('E30').Value2 = '  +4.60%  '
$ws.Range('B31').Value2 = 'Fetch.AI'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '1.45'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value2 = '  +8.31%  '
$ws.Range('B32').Value2 = 'PancakeSwap'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '2.15'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value2 = '  +1.59%  '
$ws.Range('E33').Value2 = '  -0.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '23.80'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value2 = '  +3.29%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value2 = '32.59'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value2 = '  +26.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '5.35'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value2 = '  +8.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '7.13'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value2 = '  +4.28%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value2 = '171.65'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value2 = '  +2.76%  '
$ws.Range('E39').Value2 = '  +8.73%  '
$ws.Range('D40').Value2 = '3.523.98'
$ws.Range('E41').Value2 = '  +1.76%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '0.808'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value2 = '  +4.44%  '
$ws.Range('E43').Value2 = '  +3.56%  '
$ws.Range('B44').Value2 = 'OKB'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '42.47'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value2 = '  +0.51%  '
$ws.Range('B45').Value2 = 'Stacks'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '1.73'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value2 = '  +5.66%  '
$ws.Range('E46').Value2 = '  +8.96%  '
$ws.Range('D47').Value2 = '2.624.66'
$ws.Range('E47').Value2 = '  +7.02%  '
$ws.Range('B48').Value2 = 'dogwifhat'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '2.31'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value2 = '  +15.53%  '
$ws.Range('B49').Value2 = 'InjectiveProtocol'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '23.70'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value2 = '  +7.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value2 = '6.75'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value2 = '  +1.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '0.0271'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value2 = '  +4.91%  '
